$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Date placeholder fields: 2024-02-15 -> 2024-03-12
#    (Handout Master "Date Placeholder" + Slide Master "Date Placeholder")
# ---------------------------------------------------------------------------
$hm = $p.HandoutMaster
for ($i = 1; $i -le $hm.Shapes.Count; $i++) {
    $sh = $hm.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq "2024-02-15") {
            $sh.TextFrame.TextRange.Text = "2024-03-12"
        }
    }
}

$sm = $p.SlideMaster
for ($i = 1; $i -le $sm.Shapes.Count; $i++) {
    $sh = $sm.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq "2024-02-15") {
            $sh.TextFrame.TextRange.Text = "2024-03-12"
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Slide 1 ("앞면"): 이름 box font HY엽서L -> 궁서체
#                      캠퍼스 역할 box: widen box, bold, 조선신명조 -> 바탕
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
for ($i = 1; $i -le $s1.Shapes.Count; $i++) {
    $sh = $s1.Shapes.Item($i)
    if (-not $sh.HasTextFrame) { continue }
    $txt = $sh.TextFrame.TextRange.Text
    if ($txt -eq "이름") {
        $tr = $sh.TextFrame.TextRange
        $tr.Font.Name = "궁서체"
        $tr.Font.NameFarEast = "궁서체"
    } elseif ($txt -eq "캠퍼스 역할") {
        $sh.Width = 100.7494
        $tr = $sh.TextFrame.TextRange
        $tr.Font.Bold = $true
        $tr.Font.Name = "바탕"
        $tr.Font.NameFarEast = "바탕"
    }
}

# ---------------------------------------------------------------------------
# 3) Slide 2 ("뒷면"): 캠퍼스 역할 box: widen box, bold, 조선신명조 -> 바탕
#                      이름 box font 조선신명조 -> 궁서체
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
for ($i = 1; $i -le $s2.Shapes.Count; $i++) {
    $sh = $s2.Shapes.Item($i)
    if (-not $sh.HasTextFrame) { continue }
    $txt = $sh.TextFrame.TextRange.Text
    if ($txt -eq "캠퍼스 역할") {
        $sh.Width = 100.7494
        $tr = $sh.TextFrame.TextRange
        $tr.Font.Bold = $true
        $tr.Font.Name = "바탕"
        $tr.Font.NameFarEast = "바탕"
    } elseif ($txt -eq "이름") {
        $tr = $sh.TextFrame.TextRange
        $tr.Font.Name = "궁서체"
        $tr.Font.NameFarEast = "궁서체"
    }
}
